# ajuste de alguns erros - Laura
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: mark Status ("Resolvido") for the "A Tabela 1 não é citada no texto!" item
$ws.Range("B27").Value = "Resolvido"

# Row 37: mark Status ("Resolvido") and add a note in column D explaining the fix
$ws.Range("B37").Value = "Resolvido"
$ws.Range("D37").Value = "Retirei o parágrafo"

# Row 45: new pending item with responsible person set
$ws.Range("A45").Value = "Adicionar a etapa de teste e validação no diagrama da Figura 3"
$ws.Range("C45").Value = "Laura"

# Move the visible viewport / selection to reflect where work left off
$ws.Range("C46").Select()
$excel.ActiveWindow.ScrollRow = 42
